# Dashboard_2026.xlsx budget update
# - Refresh "Atualizado" timestamp on Dashboard
# - Adjust monthly category budgets (and their dependent % / variance cells)
#   on the Dashboard sheet based on 6-month real spending data
# - Recompute the dependent variance/ratio cells on the Categorias sheet
# - Refresh the sync_timestamp on the Dados sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Dashboard sheet
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

# Last-updated banner
$dashboard.Range("A2").Value = "Atualizado: 30/12/2025 21:59"

# RESUMO DO MES - "Gastos Variaveis" projected budget (sum of the category
# budgets below) and its variance % vs. the real spend (31101.39)
$dashboard.Range("B7").Value = 17800

# D7 stores its variance as literal text ("74%"), not a real number - assigning
# a percent-looking string straight into a General-formatted cell makes Excel
# "smart" parse it into a numeric 0.74 with a new percentage number format
# (and a brand new style record). Stage the text in a scratch cell that is
# pre-formatted as Text, then copy/paste-special just the value onto D7 so the
# destination cell's existing style/format is left completely untouched.
$dashboard.Range("ZZ1").NumberFormat = "@"
$dashboard.Range("ZZ1").Value = "74%"
$dashboard.Range("ZZ1").Copy()
$dashboard.Range("D7").PasteSpecial(-4163)
$dashboard.Range("ZZ1").Clear()

# GASTOS POR CATEGORIA - per-category budgets and % of budget spent
$dashboard.Range("B14").Value = 3500   # Alimentacao
$dashboard.Range("D14").Value = 144

$dashboard.Range("B15").Value = 2500   # Transporte
$dashboard.Range("D15").Value = 315

$dashboard.Range("B16").Value = 500    # Saude
$dashboard.Range("D16").Value = 1900

$dashboard.Range("B17").Value = 4000   # Assinaturas
$dashboard.Range("D17").Value = 87

$dashboard.Range("B18").Value = 2000   # Compras
$dashboard.Range("D18").Value = 98

$dashboard.Range("B19").Value = 3500   # Lazer
$dashboard.Range("D19").Value = 72

$dashboard.Range("B21").Value = 200    # Casa
$dashboard.Range("D21").Value = 19

$dashboard.Range("B22").Value = 100    # Taxas

# ---------------------------------------------------------------------------
# Categorias sheet - variance (budget - real) and ratio (real / budget)
# recomputed against the new budgets above
# ---------------------------------------------------------------------------
$categorias = $wb.Worksheets.Item("Categorias")

$categorias.Range("D4").Value = -1548.87
$categorias.Range("E4").Value = 1.442534285714286

$categorias.Range("D5").Value = -5382.21
$categorias.Range("E5").Value = 3.152884

$categorias.Range("D6").Value = -9000
$categorias.Range("E6").Value = 19

$categorias.Range("D7").Value = 507.0900000000001
$categorias.Range("E7").Value = 0.8732274999999999

$categorias.Range("D8").Value = 32.59999999999991
$categorias.Range("E8").Value = 0.9837

$categorias.Range("D9").Value = 948.7599999999998
$categorias.Range("E9").Value = 0.7289257142857144

$categorias.Range("D11").Value = 160.1
$categorias.Range("E11").Value = 0.1995

$categorias.Range("D12").Value = 100

# ---------------------------------------------------------------------------
# Dados sheet - sync timestamp
# ---------------------------------------------------------------------------
$dados = $wb.Worksheets.Item("Dados")
$dados.Range("B3").Value = "2025-12-30T21:59:11.964253"
